$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: "Save" header, styled like the other header cells (copy format from G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Save column values (H2:H12) - 1 for the "saved" games, else 0
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
